# Update the "Terms" sheet (medical_terms.xlsx) with the latest LLM-tuning
# data: a new mis-parsed long term lands in row 2, every previously-existing
# term shifts down a row with refreshed usagecount/simpScore stats, and three
# brand-new terms ("infarction", "myocardial infarction", "myocardial") are
# appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("pneumonoultramicroscopicsilicovolcanoconiosis", 2, 0.5, 1, 0, "2025-10-23T07:41:30.284Z"),
    @("Hypertension", 8, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("Vital signs", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("BP", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("HR", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Temp", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Diagnosis", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Plan", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("mg", 4, 3.25, 0, 1, "2025-10-23T07:46:08.649Z"),
    @("daily", 4, 3.25, 0, 1, "2025-10-23T07:46:08.649Z"),
    @("Follow-up", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("atrial fibrillation", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("blood pressure", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("Temp 36", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Aspirin 81", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Palpitations", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("heartbeats", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("Examination", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("Extremities", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("Assessment", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("fibrillation", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("uncontrolled", 4, 3.25, 1, 0, "2025-10-23T07:46:08.649Z"),
    @("Peripheral", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("prevention", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("cardiology", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("evaluation", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("medication", 4, 3.25, 0, 0, "2025-10-23T07:46:08.649Z"),
    @("infarction", 2, 3, 0, 0, "2025-10-23T07:47:24.691Z"),
    @("myocardial infarction", 2, 3, 1, 0, "2025-10-23T07:47:24.693Z"),
    @("myocardial", 2, 3, 0, 0, "2025-10-23T07:47:24.693Z")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
